$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Link" header in F2, matching the formatting of the existing header row
$ws.Range("F2").Value = "Link"
$ws.Range("E2").Copy()
$ws.Range("F2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Set column F width to match the diff (stored width 38.42578125 characters)
$ws.Columns.Item(6).ColumnWidth = 37.6666666666667

# Update the selected/active cell to D5
$ws.Range("D5").Select()
